$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 46081 = 2026-02-28)
# that was bumped by one day (serial 46082 = 2026-03-01) for every data
# row (rows 2 through 246).
$lastRow = 246
$ws.Range("C2:C$lastRow").Value = 46082
